$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New tracked-date columns (row 2 headers) -----------------------------
# Three more working days are being logged: 10/29, 10/31 and 11/1/2019
# (serials 43767, 43769, 43770). The brand-new last entry (H2) keeps the
# sheet's original "d-mmm" look, while the rest of the date row moves to
# the plain short-date format.
$ws.Range("H2").Value = 43770
$ws.Range("H2").NumberFormat = "d-mmm"

$ws.Range("F2").Value = 43767
$ws.Range("G2").Value = 43769

$ws.Range("D2").NumberFormat = "mm-dd-yy"
$ws.Range("E2").NumberFormat = "mm-dd-yy"
$ws.Range("F2").NumberFormat = "mm-dd-yy"
$ws.Range("G2").NumberFormat = "mm-dd-yy"

# --- New effort-hour facts logged against those dates ---------------------
$ws.Range("F3").Value = 0.5    # Purpose, scope, definitions
$ws.Range("F8").Value = 1      # Functional requirements
$ws.Range("F10").Value = 0.5   # Formal analysis using Alloy
$ws.Range("G10").Value = 2     # Formal analysis using Alloy
$ws.Range("H10").Value = 3     # Formal analysis using Alloy

# --- Leave the selection where the author last left it --------------------
$ws.Range("H11").Select()
